$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename gun names
$ws.Range("B1").Value = "Strong"
$ws.Range("C1").Value = "Weak"

# Replace "Time to fire (frames)" row with "Rate of Fire (rounds per minute)" row
$ws.Range("A2").Value = "Rate of Fire (rounds per minute)"
$ws.Range("B2").Value = 155
$ws.Range("C2").Value = 150

# Clear the old note cell
$ws.Range("E2").Value = ""

# Adjust column A width to fit new content
$ws.Range("A1").EntireColumn.ColumnWidth = 30.7

# Update selection to reflect new active cell
$ws.Range("D6").Select()
